$d = $word.ActiveDocument

# "driverStation SMALLINT," paragraph -> apply strikethrough to whole paragraph
$rng1 = $d.Content
$null = $rng1.Find.Execute("driverStation", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para1 = $rng1.Paragraphs(1)
$para1.Range.Font.StrikeThrough = 1

# "startPlace SMALLINT, " paragraph -> apply strikethrough to whole paragraph
$rng2 = $d.Content
$null = $rng2.Find.Execute("startPlace", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$para2 = $rng2.Paragraphs(1)
$para2.Range.Font.StrikeThrough = 1
